$d = $word.ActiveDocument
$xmlNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the paragraph that holds "heloo" (it gets capitalized to "Heloo"
# and split into two runs: "H" and "eloo").
$target = $d.Paragraphs.Item(1)
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "heloo*") {
        $target = $p
    }
}

# Rebuild the paragraph's contents as two runs: "H" + "eloo" (this also
# capitalizes the leading "h"). InsertXML on the text-only portion of the
# paragraph (excluding the paragraph mark) swaps in a clean two-run split
# with no leftover empty rPr on either run (unlike toggling a format on
# then off), while leaving the paragraph's own identity/properties alone.
$textRange = $d.Range($target.Range.Start, $target.Range.End - 1)
$newRunsXml = "<w:p $xmlNs><w:r><w:t>H</w:t></w:r><w:r><w:t>eloo</w:t></w:r></w:p>"
$textRange.InsertXML($newRunsXml) | Out-Null

# Re-fetch that paragraph (it is still the first one) and add a new
# paragraph right after it containing the word "fourth".
$target = $d.Paragraphs.Item(1)
$afterPoint = $target.Range.End
$insertionRange = $d.Range($afterPoint, $afterPoint)
$insertionRange.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item(2)
$newPara.Range.InsertAfter("fourth") | Out-Null
